# feat: add 2022-Q1 data
#
# 1) A new sheet "2022-Q1" is inserted right before the "总计" (totals) sheet,
#    holding the per-fund holdings snapshot for 2022-Q1 (same shape as the
#    other quarterly sheets).
# 2) The "总计" (totals) sheet gets a new first data row for "2022-Q1"
#    (4 funds held, 2.36 billion yuan), with the existing rows shifting down
#    and the running index in column A renumbered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" worksheet by cloning the "2021-Q4" sheet's
#    layout (same headers/styles) and place it immediately before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item("总计"))

# The copy is inserted right before "总计" and named "2021-Q4 (2)".
$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# The template sheet had 7 data rows (A2:H8); the 2022-Q1 snapshot only has
# 4, so drop the trailing rows entirely (not just their contents) so the
# sheet's dimension shrinks back down to A1:H5.
$ws.Range("A6:H8").Delete()

# Force the fund-code / numeric-looking text columns to be stored as text
# so values like "008099" keep their leading zeros instead of becoming 8099.
$ws.Range("B2:G5").NumberFormat = "@"

$fundRows = @(
  @("008099", "广发价值领先混合",               "61.82", "83.88", "3.79", "2.3430", 8),
  @("003749", "创金合信鑫收益灵活配置混合A",     "0.65",  "51.22", "1.30", "0.0084", 8),
  @("006906", "创金合信鑫收益灵活配置混合E",     "0.65",  "51.22", "1.30", "0.0084", 8),
  @("003750", "创金合信鑫收益灵活配置混合C",     "0.02",  "51.22", "1.30", "0.0003", 8)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $i
  $ws.Cells.Item($row, 2).Value = $fundRows[$i][0]
  $ws.Cells.Item($row, 3).Value = $fundRows[$i][1]
  $ws.Cells.Item($row, 4).Value = $fundRows[$i][2]
  $ws.Cells.Item($row, 5).Value = $fundRows[$i][3]
  $ws.Cells.Item($row, 6).Value = $fundRows[$i][4]
  $ws.Cells.Item($row, 7).Value = $fundRows[$i][5]
  $ws.Cells.Item($row, 8).Value = $fundRows[$i][6]
}

# The "@" number format above served only to keep the text values from
# being re-interpreted as numbers; drop it again afterwards so these
# cells end up with the same (unstyled) look as the other quarter sheets.
$ws.Range("B2:G5").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new top data row for 2022-Q1 and
#    renumber the running index in column A for the rows pushed down.
# Re-fetch the handle by name now (sheet handles captured before the
# insert above point at a *position*, not a stable sheet identity, and
# would now resolve to the wrong worksheet).
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$zongji.Rows(2).Insert()

# The inserted row inherits odd formatting from the row above it; clear
# that back to the unstyled look the other B:D data cells use, then copy
# just the index-column style (bold/centered/bordered) from row 3.
$zongji.Range("B2:D2").ClearFormats()
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 4
$zongji.Range("D2").Value = 2.36

for ($r = 3; $r -le 7; $r++) {
  $zongji.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "done"
